# Add "Gameweeks import" support columns to the Challenges sheet.
$wb = $excel.ActiveWorkbook

# Work explicitly against the "Challenges" sheet (the active one).
$ws = $wb.Worksheets.Item("Challenges")

# New header columns, appended after the existing "Skill LongTermVision" (R) column.
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# Sample/default data row values for the new columns.
# Leading apostrophe forces literal text "true" (not an auto-converted boolean);
# resetting the style afterwards drops the quote-prefix formatting it implies.
$flag = $ws.Range("S2")
$flag.Value = "'true"
$flag.Style = "Normal"
$ws.Range("T2").Value = 1
